# Update gh-pages to output generated at 456a3b4
#
# Changes applied (per the canonical OOXML diff):
#  - "展览" (sheet 1) and "全部类型" (sheet 4) sheets: several "想去人数"
#    (F column) counters bumped up by a handful, and a brand-new event
#    ("苏州·星部落动漫嘉年华") inserted right before the trailing two rows,
#    pushing the previously-last two rows down by one.
#  - "演出" (sheet 2): one F-column counter bumped by 1.
#  - "本地生活" (sheet 3): untouched.

$wb = $excel.ActiveWorkbook

function Bump-F {
    param($ws, $row, $newValue)
    $ws.Cells.Item($row, 6).Value = $newValue
}

function Insert-NewFirstOfLastThree {
    param($ws, $insertAt)

    # $insertAt is the row that currently holds the (soon to be second-to-
    # last) "Redamancy" event. Opening a blank row here pushes it (and the
    # still-last "理想乡" row after it) down by one, then the new
    # "星部落" event takes over the now-vacated row $insertAt.
    $ws.Cells.Item($insertAt, 1).EntireRow.Insert()

    # Row-insert sometimes drops the row-index column's border when it
    # blends formatting from its neighbours; re-stamp it from the row right
    # below (which still carries the original, untouched formatting).
    $srcRow = $insertAt + 1
    $ws.Range("A$srcRow`:I$srcRow").Copy() | Out-Null
    $ws.Range("A$insertAt`:I$insertAt").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    # The new row keeps the same sequential index the old row had (A = row-1).
    $seq = $insertAt - 1

    $ws.Cells.Item($insertAt, 1).Value = $seq

    # Force the date column to stay plain text (otherwise Excel helpfully
    # reinterprets an ISO-looking string as a date serial number).
    $dateCell = $ws.Cells.Item($insertAt, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2024-08-03"
    $dateCell.ClearFormats()

    $ws.Cells.Item($insertAt, 3).Value = "苏州·星部落动漫嘉年华"
    $ws.Cells.Item($insertAt, 4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
    $ws.Cells.Item($insertAt, 5).Value = "2024.08.03 09:00-08.04 16:00"
    $ws.Cells.Item($insertAt, 6).Value = 2
    $ws.Cells.Item($insertAt, 7).Value = 49
    $ws.Cells.Item($insertAt, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84858"
    $ws.Cells.Item($insertAt, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"

    # The two rows that got shifted down keep their data, but their
    # sequential index (column A) needs to be bumped by one to stay
    # consistent with their new row number.
    $ws.Cells.Item($insertAt + 1, 1).Value = $seq + 1
    $ws.Cells.Item($insertAt + 2, 1).Value = $seq + 2
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Bump-F $ws1 2  14951
Bump-F $ws1 3  18798
Bump-F $ws1 5  132
Bump-F $ws1 14 131
Bump-F $ws1 20 92
Bump-F $ws1 22 7812
Bump-F $ws1 24 30
Bump-F $ws1 26 1233
Bump-F $ws1 28 6002
Bump-F $ws1 34 5373

Insert-NewFirstOfLastThree $ws1 35

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

Bump-F $ws2 3 12

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 -- no changes
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Bump-F $ws4 2  14951
Bump-F $ws4 3  18798
Bump-F $ws4 5  132
Bump-F $ws4 14 131
Bump-F $ws4 21 92
Bump-F $ws4 23 7812
Bump-F $ws4 25 30
Bump-F $ws4 27 1233
Bump-F $ws4 29 12
Bump-F $ws4 31 6002
Bump-F $ws4 37 5373

Insert-NewFirstOfLastThree $ws4 38
